$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Reorder the "estado de cuenta" data rows (16-22):
# The GINA PATRICIA BAHOQUE LOPEZ periods (1608..1612) now occupy rows 17-21
# in ascending order, and LUIS MANUEL ANAYA PASO / 1705 moves down to row 22.
# Row 16 (VICTOR ALFONSO ESCALANTE LARA / 1607) is unchanged.

$ws.Range("C17").Value = "32907894"
$ws.Range("D17").Value = "GINA PATRICIA BAHOQUE LOPEZ"
$ws.Range("E17").Value = "1608"
$ws.Range("F17").Value = 27578
$ws.Range("G17").Value = 737717

$ws.Range("C18").Value = "32907894"
$ws.Range("D18").Value = "GINA PATRICIA BAHOQUE LOPEZ"
$ws.Range("E18").Value = "1609"
$ws.Range("F18").Value = 27578
$ws.Range("G18").Value = 737717

$ws.Range("C19").Value = "32907894"
$ws.Range("D19").Value = "GINA PATRICIA BAHOQUE LOPEZ"
$ws.Range("E19").Value = "1610"
$ws.Range("F19").Value = 27578
$ws.Range("G19").Value = 737717

$ws.Range("C20").Value = "32907894"
$ws.Range("D20").Value = "GINA PATRICIA BAHOQUE LOPEZ"
$ws.Range("E20").Value = "1611"
$ws.Range("F20").Value = 27578
$ws.Range("G20").Value = 737717

$ws.Range("C21").Value = "32907894"
$ws.Range("D21").Value = "GINA PATRICIA BAHOQUE LOPEZ"
$ws.Range("E21").Value = "1612"
$ws.Range("F21").Value = 27578
$ws.Range("G21").Value = 737717

$ws.Range("C22").Value = "1143329524"
$ws.Range("D22").Value = "LUIS MANUEL ANAYA PASO"
$ws.Range("E22").Value = "1705"
$ws.Range("F22").Value = 29509
$ws.Range("G22").Value = 737717
